# Update "PoFDCtAE" sheet: replace formulas referencing 'Data from BFPIaE'
# with literal value 1 for a set of cells (imported elec/fuels update).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PoFDCtAE")

$cellsToOne = @("C3", "D4", "I9", "J10", "K11", "L12", "M13", "N14", "S19", "T20")
foreach ($ref in $cellsToOne) {
    $ws.Range($ref).Value = 1
}

# Recalculate so dependent formulas (R10, R11, R14, R19, R20 = 1 - <cell>)
# pick up the new values.
$excel.Calculate()

# Restore the active selection on this sheet to match the saved workbook
# state (cell U20 selected in the bottom-right frozen pane), without
# leaving this sheet as the active tab.
$originalActive = $wb.ActiveSheet
$ws.Range("U20").Select()
$originalActive.Activate()
